$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates
# Cells whose new value is unambiguously numeric-looking need NumberFormat = "@"
# applied first, so Excel keeps them as text (matching the original inlineStr data)
# instead of silently converting them to numbers.
$dTextSafeValues = @{
    2 = "34.503.05"
    3 = "1.810.06"
    12 = "2.071.69"
    14 = "1.814.13"
    16 = "34.465.16"
    20 = "0.0₃0775"
    35 = "1.363.91"
    47 = "1.971.45"
}

$dNumericLookingValues = @{
    5 = "226.05"
    8 = "36.32"
    13 = "11.29"
    18 = "68.65"
    19 = "242.79"
    25 = "171.53"
    26 = "7.92"
    27 = "17.29"
    33 = "0.0518"
    38 = "2.36"
    40 = "2.42"
    42 = "81.01"
    43 = "0.937"
    45 = "13.37"
    46 = "0.0499"
    50 = "102.82"
}

foreach ($row in $dTextSafeValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $dTextSafeValues[$row]
}

foreach ($row in $dNumericLookingValues.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $dNumericLookingValues[$row]
}

# Volume(1h) (column E) updates - these strings include leading/trailing spaces
# and a percent sign, so Excel naturally keeps them as text.
$eValues = @{
    2 = "  +0.09%  "
    3 = "  +0.26%  "
    4 = "  -0.10%  "
    5 = "  -1.01%  "
    6 = "  +2.69%  "
    7 = "  -0.09%  "
    8 = "  +3.77%  "
    9 = "  -2.25%  "
    10 = "  -2.04%  "
    11 = "  +1.49%  "
    12 = "  +0.34%  "
    13 = "  +0.37%  "
    14 = "  +0.53%  "
    15 = "  -1.92%  "
    16 = "  +0.07%  "
    17 = "  +1.41%  "
    18 = "  -0.66%  "
    19 = "  -1.26%  "
    20 = "  -2.91%  "
    21 = "  -2.57%  "
    22 = "  -0.11%  "
    23 = "  -1.45%  "
    24 = "  +4.91%  "
    25 = "  -1.22%  "
    26 = "  +2.89%  "
    27 = "  +2.98%  "
    28 = "  +1.91%  "
    29 = "  -0.13%  "
    30 = "  -0.15%  "
    31 = "  -1.92%  "
    32 = "  -0.96%  "
    33 = "  -2.04%  "
    34 = "  -1.99%  "
    35 = "  -2.38%  "
    36 = "  -4.06%  "
    37 = "  -0.51%  "
    38 = "  -6.19%  "
    39 = "  -1.90%  "
    40 = "  +1.76%  "
    41 = "  -1.58%  "
    42 = "  -2.84%  "
    43 = "  -1.46%  "
    44 = "  +4.40%  "
    45 = "  -1.44%  "
    46 = "  -2.10%  "
    47 = "  +0.33%  "
    48 = "  -2.51%  "
    49 = "  -0.14%  "
    50 = "  -2.01%  "
    51 = "  -5.12%  "
}

foreach ($row in $eValues.Keys) {
    $ws.Cells.Item($row, 5).Value = $eValues[$row]
}
